# Update "想去人数" (F column) figures across all sheets, and fix one
# "最低票价" (G column) cell that changed from inline text "不可售" to a
# numeric value, matching the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 439
$ws.Range("F3").Value = 2765
$ws.Range("F4").Value = 1325
$ws.Range("F5").Value = 60
$ws.Range("F7").Value = 583
$ws.Range("F9").Value = 598
$ws.Range("F10").Value = 277
$ws.Range("F11").Value = 84
$ws.Range("F12").Value = 11511
$ws.Range("F13").Value = 6549
$ws.Range("F19").Value = 11
$ws.Range("F20").Value = 912
$ws.Range("F21").Value = 59
$ws.Range("F22").Value = 261
$ws.Range("F23").Value = 918
$ws.Range("F24").Value = 3632
$ws.Range("F27").Value = 494
$ws.Range("F28").Value = 161
$ws.Range("F29").Value = 311
$ws.Range("F30").Value = 15
$ws.Range("F31").Value = 265
$ws.Range("G33").Value = 1
$ws.Range("F35").Value = 1228
$ws.Range("F36").Value = 224
$ws.Range("F37").Value = 418
$ws.Range("F38").Value = 189
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 36
$ws.Range("F6").Value = 24
$ws.Range("F8").Value = 139
$ws.Range("F11").Value = 3669
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 79
$ws.Range("F23").Value = 14
$ws.Range("F25").Value = 8
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9008
$ws.Range("F3").Value = 492
$ws.Range("F4").Value = 1806
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9008
$ws.Range("F3").Value = 492
$ws.Range("F4").Value = 1806
$ws.Range("F5").Value = 439
$ws.Range("F6").Value = 2765
$ws.Range("F8").Value = 36
$ws.Range("F9").Value = 24
$ws.Range("F10").Value = 1325
$ws.Range("F12").Value = 60
$ws.Range("F13").Value = 583
$ws.Range("F15").Value = 139
$ws.Range("F16").Value = 598
$ws.Range("F17").Value = 277
$ws.Range("F18").Value = 84
$ws.Range("F19").Value = 11511
$ws.Range("F20").Value = 3669
$ws.Range("F27").Value = 11
$ws.Range("F28").Value = 912
$ws.Range("F29").Value = 59
$ws.Range("F30").Value = 261
$ws.Range("F31").Value = 918
$ws.Range("F32").Value = 3632
$ws.Range("F35").Value = 161
$ws.Range("F36").Value = 311
$ws.Range("F37").Value = 265
$ws.Range("F41").Value = 1228
$ws.Range("F42").Value = 224
$ws.Range("F43").Value = 189
$ws.Range("F46").Value = 79
$ws.Range("F47").Value = 14
$ws.Range("F49").Value = 8

Write-Output "Applied 68 cell updates"
